$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that currently sits on its own
#    empty paragraph (right before the "Secinajumi" list item).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Locate the "Pielikums. Nokopet kodu ..." paragraph (last
#    non-empty paragraph) and append a trailing space run, followed
#    by a fresh "_GoBack" bookmark at the very end of the paragraph.
# ------------------------------------------------------------------
$total = $d.Paragraphs.Count
$target = $null
for ($i = 1; $i -le $total; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Nokop*kodu*") {
        $target = $cand
    }
}

$r = $target.Range
$body = $d.Range($r.Start, $r.End - 1)
$body.InsertAfter(" ")

# Re-resolve the paragraph range after the mutation above.
$r2 = $target.Range
$body2 = $d.Range($r2.Start, $r2.End - 1)

# The runtime mis-resolves a *zero length* Range sitting exactly on a
# paragraph-end boundary, so anchor the bookmark using a temporary
# marker character placed right after the space, then remove the
# marker once the bookmark has been created at the safe position
# between the space and the marker.
$marker = [char]1
$body2.InsertAfter([string]$marker)
$r3 = $target.Range
$anchorPos = $r3.End - 2
$bmRange = $d.Range($anchorPos, $anchorPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Range($anchorPos, $anchorPos + 1)
$markerRange.Delete()

# ------------------------------------------------------------------
# 3) Add a brand-new list paragraph right after it containing the
#    GitHub link, inheriting the numbering / character formatting of
#    the paragraph above.
# ------------------------------------------------------------------
$r4 = $target.Range
$body4 = $d.Range($r4.Start, $r4.End - 1)
$body4.InsertParagraphAfter()

$newPara = $target.Next()
$nr = $newPara.Range
$nbody = $d.Range($nr.Start, $nr.End - 1)
$nbody.InsertAfter("https://github.com/PatrUpe600/08_04_Patriks")
